# schedule.xlsx update: "Updated slides, shedule, notes"
#
# The underlying change re-purposes the old "Scales" topic row (15) into the
# "Summarizing data" topic (pulling in its slides/notes links), re-purposes
# the old "Summarizing data" topic row (16) into the "Plotting facets" topic
# (which previously lived as its own placeholder row 24, now removed), and
# adds video links for both. Row 17 ("Reading data") keeps its own content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: was "Scales" -> becomes "Summarizing data" ---
$ws.Range("D15").Value = "Summarizing data"
$ws.Range("G15").Value = "Video"
$ws.Range("L15").Value = "slides/10-summarizing-data.html"

# --- Row 16: was "Summarizing data" -> becomes "Plotting facets" ---
$ws.Range("D16").Value = "Plotting facets"
$ws.Range("G16").Value = "Video"
$ws.Range("J16").Value = "#facets"
$ws.Range("L16").Value = "slides/09-facets.html"
$ws.Range("K16").Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=b9f18ad7-21a3-4c07-a912-acab001bb44f"

# --- back to row 15's remaining link cells ---
$ws.Range("J15").Value = "#summarizing-data"
$ws.Range("K15").Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=f751632b-e900-4121-80b8-acab01616c37"

# --- Row 24: drop the now-redundant standalone "Plotting facets" entry ---
$ws.Range("D24").ClearContents()

# --- Selection moves from K16 to K15 ---
$ws.Range("K15").Select()
